$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (naive -> joint-datasets)
$ws.Name = "joint-datasets_after_task1"

# Update existing rows (3-48) with refreshed training metrics
# Row 3
$ws.Cells.Item(3, 3).Value = 3.918494273291694
$ws.Cells.Item(3, 4).Value = 3.668997688293457
$ws.Cells.Item(3, 5).Value = 13.71
$ws.Cells.Item(3, 9).Value = 0.01824049613475799
$ws.Cells.Item(3, 10).Value = 13.71

# Row 4
$ws.Cells.Item(4, 3).Value = 3.329281957414415
$ws.Cells.Item(4, 4).Value = 3.310175447463989
$ws.Cells.Item(4, 5).Value = 20.09
$ws.Cells.Item(4, 9).Value = 0.01634105877876282
$ws.Cells.Item(4, 10).Value = 20.09

# Row 5
$ws.Cells.Item(5, 3).Value = 2.956839343176948
$ws.Cells.Item(5, 4).Value = 3.040581960678101
$ws.Cells.Item(5, 5).Value = 24.24
$ws.Cells.Item(5, 9).Value = 0.01497891943454742
$ws.Cells.Item(5, 10).Value = 24.24

# Row 6
$ws.Cells.Item(6, 3).Value = 2.674733493593004
$ws.Cells.Item(6, 4).Value = 2.768585596084595
$ws.Cells.Item(6, 5).Value = 29.61
$ws.Cells.Item(6, 9).Value = 0.01366797437667847
$ws.Cells.Item(6, 10).Value = 29.61

# Row 7
$ws.Cells.Item(7, 3).Value = 2.452956551445855
$ws.Cells.Item(7, 4).Value = 2.84297324180603
$ws.Cells.Item(7, 5).Value = 27.97
$ws.Cells.Item(7, 9).Value = 0.01415626933574676
$ws.Cells.Item(7, 10).Value = 27.97

# Row 8
$ws.Cells.Item(8, 3).Value = 2.273474664688111
$ws.Cells.Item(8, 4).Value = 2.474572916030884
$ws.Cells.Item(8, 5).Value = 36.25
$ws.Cells.Item(8, 9).Value = 0.01214034621715546
$ws.Cells.Item(8, 10).Value = 36.25

# Row 9
$ws.Cells.Item(9, 3).Value = 2.130468515290154
$ws.Cells.Item(9, 4).Value = 2.62041241645813
$ws.Cells.Item(9, 5).Value = 33.3
$ws.Cells.Item(9, 9).Value = 0.01307356290817261
$ws.Cells.Item(9, 10).Value = 33.3

# Row 10
$ws.Cells.Item(10, 3).Value = 2.007545682059394
$ws.Cells.Item(10, 4).Value = 2.417155933380127
$ws.Cells.Item(10, 5).Value = 38.3
$ws.Cells.Item(10, 9).Value = 0.01200385262966156
$ws.Cells.Item(10, 10).Value = 38.3

# Row 11
$ws.Cells.Item(11, 3).Value = 1.89463817331526
$ws.Cells.Item(11, 4).Value = 2.238714647293091
$ws.Cells.Item(11, 5).Value = 41.03
$ws.Cells.Item(11, 9).Value = 0.0111312420129776
$ws.Cells.Item(11, 10).Value = 41.03

# Row 12
$ws.Cells.Item(12, 3).Value = 1.800584742758009
$ws.Cells.Item(12, 4).Value = 2.225303735733032
$ws.Cells.Item(12, 5).Value = 41.01
$ws.Cells.Item(12, 9).Value = 0.01114006059169769
$ws.Cells.Item(12, 10).Value = 41.01

# Row 13
$ws.Cells.Item(13, 3).Value = 1.710838945176866
$ws.Cells.Item(13, 4).Value = 2.150233263969421
$ws.Cells.Item(13, 5).Value = 43.1
$ws.Cells.Item(13, 9).Value = 0.0106656606554985
$ws.Cells.Item(13, 10).Value = 43.1

# Row 14
$ws.Cells.Item(14, 3).Value = 1.633768646452162
$ws.Cells.Item(14, 4).Value = 2.111660628318786
$ws.Cells.Item(14, 5).Value = 43.15
$ws.Cells.Item(14, 9).Value = 0.01065636731386185
$ws.Cells.Item(14, 10).Value = 43.15

# Row 15
$ws.Cells.Item(15, 3).Value = 1.553494999143812
$ws.Cells.Item(15, 4).Value = 2.383521327972412
$ws.Cells.Item(15, 5).Value = 40.56
$ws.Cells.Item(15, 9).Value = 0.01186143044233322
$ws.Cells.Item(15, 10).Value = 40.56

# Row 16
$ws.Cells.Item(16, 3).Value = 1.482246497472127
$ws.Cells.Item(16, 4).Value = 2.326114540100098
$ws.Cells.Item(16, 5).Value = 41.99
$ws.Cells.Item(16, 9).Value = 0.0115641282081604
$ws.Cells.Item(16, 10).Value = 41.99

# Row 17
$ws.Cells.Item(17, 3).Value = 1.417229968706767
$ws.Cells.Item(17, 4).Value = 2.259857697486877
$ws.Cells.Item(17, 5).Value = 41.86
$ws.Cells.Item(17, 9).Value = 0.01131694558858871
$ws.Cells.Item(17, 10).Value = 41.86

# Row 18
$ws.Cells.Item(18, 3).Value = 1.342169585757785
$ws.Cells.Item(18, 4).Value = 2.154725122451782
$ws.Cells.Item(18, 5).Value = 44.62
$ws.Cells.Item(18, 9).Value = 0.01073354250192642
$ws.Cells.Item(18, 10).Value = 44.62

# Row 19
$ws.Cells.Item(19, 3).Value = 1.283520460658603
$ws.Cells.Item(19, 4).Value = 2.182112832069397
$ws.Cells.Item(19, 5).Value = 45.15
$ws.Cells.Item(19, 9).Value = 0.01083021525144577
$ws.Cells.Item(19, 10).Value = 45.15

# Row 20
$ws.Cells.Item(20, 3).Value = 1.356951496866014
$ws.Cells.Item(20, 4).Value = 1.902456116676331
$ws.Cells.Item(20, 5).Value = 48.55
$ws.Cells.Item(20, 9).Value = 0.009425257515907287
$ws.Cells.Item(20, 10).Value = 48.55

# Row 21
$ws.Cells.Item(21, 3).Value = 1.28929979801178
$ws.Cells.Item(21, 4).Value = 1.909884910583496
$ws.Cells.Item(21, 5).Value = 48.75
$ws.Cells.Item(21, 9).Value = 0.009460268962383269
$ws.Cells.Item(21, 10).Value = 48.75

# Row 22
$ws.Cells.Item(22, 3).Value = 1.252987093925476
$ws.Cells.Item(22, 4).Value = 1.924414305686951
$ws.Cells.Item(22, 5).Value = 48.8
$ws.Cells.Item(22, 9).Value = 0.009552257227897643
$ws.Cells.Item(22, 10).Value = 48.8

# Row 23
$ws.Cells.Item(23, 3).Value = 1.218699714342753
$ws.Cells.Item(23, 4).Value = 1.96317008972168
$ws.Cells.Item(23, 5).Value = 47.8
$ws.Cells.Item(23, 9).Value = 0.009688619220256805
$ws.Cells.Item(23, 10).Value = 47.8

# Row 24
$ws.Cells.Item(24, 3).Value = 1.185857849386003
$ws.Cells.Item(24, 4).Value = 1.952775983810425
$ws.Cells.Item(24, 5).Value = 48.17
$ws.Cells.Item(24, 9).Value = 0.009651916718482971
$ws.Cells.Item(24, 10).Value = 48.17

# Row 25
$ws.Cells.Item(25, 3).Value = 1.153948659631941
$ws.Cells.Item(25, 4).Value = 1.965516324043274
$ws.Cells.Item(25, 5).Value = 47.92
$ws.Cells.Item(25, 9).Value = 0.009790390431880951
$ws.Cells.Item(25, 10).Value = 47.92

# Row 26
$ws.Cells.Item(26, 3).Value = 1.262865560319689
$ws.Cells.Item(26, 4).Value = 1.884373531341553
$ws.Cells.Item(26, 5).Value = 49.28
$ws.Cells.Item(26, 9).Value = 0.009309570682048797
$ws.Cells.Item(26, 10).Value = 49.28

# Row 27
$ws.Cells.Item(27, 3).Value = 1.247145944701301
$ws.Cells.Item(27, 4).Value = 1.883588781356812
$ws.Cells.Item(27, 5).Value = 49.21
$ws.Cells.Item(27, 9).Value = 0.009320940446853637
$ws.Cells.Item(27, 10).Value = 49.21

# Row 28
$ws.Cells.Item(28, 3).Value = 1.234420971075694
$ws.Cells.Item(28, 4).Value = 1.888651022911072
$ws.Cells.Item(28, 5).Value = 49.27
$ws.Cells.Item(28, 9).Value = 0.009336944842338562
$ws.Cells.Item(28, 10).Value = 49.27

# Row 29
$ws.Cells.Item(29, 3).Value = 1.226133259137472
$ws.Cells.Item(29, 4).Value = 1.889580140113831
$ws.Cells.Item(29, 5).Value = 49.1
$ws.Cells.Item(29, 9).Value = 0.009347826743125915
$ws.Cells.Item(29, 10).Value = 49.1

# Row 30
$ws.Cells.Item(30, 3).Value = 1.214604782528347
$ws.Cells.Item(30, 4).Value = 1.890507183074951
$ws.Cells.Item(30, 5).Value = 49.34
$ws.Cells.Item(30, 9).Value = 0.00935202819108963
$ws.Cells.Item(30, 10).Value = 49.34

# Row 31
$ws.Cells.Item(31, 3).Value = 1.207312098079258
$ws.Cells.Item(31, 4).Value = 1.89457049369812
$ws.Cells.Item(31, 5).Value = 48.99
$ws.Cells.Item(31, 9).Value = 0.009381552672386169
$ws.Cells.Item(31, 10).Value = 48.99

# Row 32
$ws.Cells.Item(32, 3).Value = 1.197794361379412
$ws.Cells.Item(32, 4).Value = 1.89603512763977
$ws.Cells.Item(32, 5).Value = 49.26
$ws.Cells.Item(32, 9).Value = 0.009388120913505554
$ws.Cells.Item(32, 10).Value = 49.26

# Row 33
$ws.Cells.Item(33, 3).Value = 1.229177800284492
$ws.Cells.Item(33, 4).Value = 1.883371248245239
$ws.Cells.Item(33, 5).Value = 49.36
$ws.Cells.Item(33, 9).Value = 0.009314729177951812
$ws.Cells.Item(33, 10).Value = 49.36

# Row 34
$ws.Cells.Item(34, 3).Value = 1.225623777177599
$ws.Cells.Item(34, 4).Value = 1.884956669807434
$ws.Cells.Item(34, 5).Value = 49.34
$ws.Cells.Item(34, 9).Value = 0.009319575476646424
$ws.Cells.Item(34, 10).Value = 49.34

# Row 35
$ws.Cells.Item(35, 3).Value = 1.225096664163801
$ws.Cells.Item(35, 4).Value = 1.883209114074707
$ws.Cells.Item(35, 5).Value = 49.27
$ws.Cells.Item(35, 9).Value = 0.009322001469135284
$ws.Cells.Item(35, 10).Value = 49.27

# Row 36
$ws.Cells.Item(36, 3).Value = 1.220825915336609
$ws.Cells.Item(36, 4).Value = 1.883412671089172
$ws.Cells.Item(36, 5).Value = 49.41
$ws.Cells.Item(36, 9).Value = 0.009321438050270081
$ws.Cells.Item(36, 10).Value = 49.41

# Row 37
$ws.Cells.Item(37, 3).Value = 1.220025650130378
$ws.Cells.Item(37, 4).Value = 1.884477229118347
$ws.Cells.Item(37, 5).Value = 49.47
$ws.Cells.Item(37, 9).Value = 0.009328422248363494
$ws.Cells.Item(37, 10).Value = 49.47

# Row 38
$ws.Cells.Item(38, 3).Value = 1.215771516693963
$ws.Cells.Item(38, 4).Value = 1.88389265537262
$ws.Cells.Item(38, 5).Value = 49.49
$ws.Cells.Item(38, 9).Value = 0.009326774680614472
$ws.Cells.Item(38, 10).Value = 49.49

# Row 39
$ws.Cells.Item(39, 3).Value = 1.214217896991306
$ws.Cells.Item(39, 4).Value = 1.887519946098328
$ws.Cells.Item(39, 5).Value = 49.48
$ws.Cells.Item(39, 9).Value = 0.009334433019161224
$ws.Cells.Item(39, 10).Value = 49.48

# Row 40
$ws.Cells.Item(40, 3).Value = 1.212401366233826
$ws.Cells.Item(40, 4).Value = 1.884832863807678
$ws.Cells.Item(40, 5).Value = 49.39
$ws.Cells.Item(40, 9).Value = 0.009330078661441803
$ws.Cells.Item(40, 10).Value = 49.39

# Row 41
$ws.Cells.Item(41, 3).Value = 1.221589766608344
$ws.Cells.Item(41, 4).Value = 1.882504839897156
$ws.Cells.Item(41, 5).Value = 49.25
$ws.Cells.Item(41, 9).Value = 0.009314165782928466
$ws.Cells.Item(41, 10).Value = 49.25

# Row 42
$ws.Cells.Item(42, 3).Value = 1.218060270945231
$ws.Cells.Item(42, 4).Value = 1.883817148208618
$ws.Cells.Item(42, 5).Value = 49.3
$ws.Cells.Item(42, 9).Value = 0.009318997454643249
$ws.Cells.Item(42, 10).Value = 49.3

# Row 43
$ws.Cells.Item(43, 3).Value = 1.219218468930986
$ws.Cells.Item(43, 4).Value = 1.884390201568604
$ws.Cells.Item(43, 5).Value = 49.4
$ws.Cells.Item(43, 9).Value = 0.009319932687282562
$ws.Cells.Item(43, 10).Value = 49.4

# Row 44
$ws.Cells.Item(44, 3).Value = 1.21764141480128
$ws.Cells.Item(44, 4).Value = 1.882481722831726
$ws.Cells.Item(44, 5).Value = 49.57
$ws.Cells.Item(44, 9).Value = 0.009319509196281434
$ws.Cells.Item(44, 10).Value = 49.57

# Row 45
$ws.Cells.Item(45, 3).Value = 1.218737951119741
$ws.Cells.Item(45, 4).Value = 1.884937329292297
$ws.Cells.Item(45, 5).Value = 49.37
$ws.Cells.Item(45, 9).Value = 0.009323397338390351
$ws.Cells.Item(45, 10).Value = 49.37

# Row 46
$ws.Cells.Item(46, 3).Value = 1.217519544760386
$ws.Cells.Item(46, 4).Value = 1.884719605445862
$ws.Cells.Item(46, 5).Value = 49.28
$ws.Cells.Item(46, 9).Value = 0.009320978903770446
$ws.Cells.Item(46, 10).Value = 49.28

# Row 47
$ws.Cells.Item(47, 3).Value = 1.21656859503852
$ws.Cells.Item(47, 4).Value = 1.884464273452759
$ws.Cells.Item(47, 5).Value = 49.38
$ws.Cells.Item(47, 9).Value = 0.00932421624660492
$ws.Cells.Item(47, 10).Value = 49.38

# Row 48
$ws.Cells.Item(48, 3).Value = 1.21674672709571
$ws.Cells.Item(48, 4).Value = 1.885056939125061
$ws.Cells.Item(48, 5).Value = 49.52
$ws.Cells.Item(48, 9).Value = 0.009326378464698791
$ws.Cells.Item(48, 10).Value = 49.52

# Append new row 49 (epoch 47)
$ws.Cells.Item(49, 1).Value = 1
$ws.Cells.Item(49, 2).Value = 47
$ws.Cells.Item(49, 3).Value = 1.217466417683496
$ws.Cells.Item(49, 4).Value = 1.883657088279724
$ws.Cells.Item(49, 5).Value = 49.46
$ws.Cells.Item(49, 7).Value = 1
$ws.Cells.Item(49, 8).Value = 47
$ws.Cells.Item(49, 9).Value = 0.00932148848772049
$ws.Cells.Item(49, 10).Value = 49.46
